# The commit updates the article Q&A sheet: the last question (about "obfit")
# was removed from the sheet and from the shared-string table, and the
# previously selected cell moved from D4 to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the question that used to live in D3 ("Cómo utilizarón la reducción
# de moleculas a través de obfit? ..."). Clear() drops both the value and the
# style, leaving no cell record behind, matching the removed <c> in the diff.
$ws.Range("D3").Clear()

# Reflect the new active selection recorded in the saved sheet view.
$ws.Range("A5").Select()
